# Publish pass: turn each section's Heading1/Heading2 paragraph into a
# Markdown-style "#"/"##" heading line styled as FirstParagraph, and demote
# the paragraph that used to follow it (styled FirstParagraph) to BodyText.
# Applies to the three sections in the document: "Clear Screen and Desk"
# (H1 -> "#..."), "Clear Desk" (H2 -> "##...") and "Feedback" (H2 -> "##...").
$d = $word.ActiveDocument

# --- Section 1: "Clear Screen and Desk" -------------------------------
# Old: Heading1 paragraph "Clear Screen and Desk" immediately followed by a
#      FirstParagraph paragraph "Users shall comply with the following:"
# New: FirstParagraph paragraph "#Clear Screen and Desk" immediately
#      followed by a BodyText paragraph "Users shall comply with the
#      following:"
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Text = "#Clear Screen and Desk"
$p1.Style = "First Paragraph"

$p2 = $d.Paragraphs.Item(2)
$p2.Style = "Body Text"

# --- Section 2: "Clear Desk" -------------------------------------------
# Old: Heading2 paragraph "Clear Desk" immediately followed by a
#      FirstParagraph paragraph "Users shall comply with the following:"
# New: FirstParagraph paragraph "##Clear Desk" immediately followed by a
#      BodyText paragraph "Users shall comply with the following:"
$p9 = $d.Paragraphs.Item(9)
$p9.Range.Text = "##Clear Desk"
$p9.Style = "First Paragraph"

$p10 = $d.Paragraphs.Item(10)
$p10.Style = "Body Text"

# --- Section 3: "Feedback" ----------------------------------------------
# Old: Heading2 paragraph "Feedback" immediately followed by a
#      FirstParagraph paragraph with the contact sentence + hyperlink.
# New: FirstParagraph paragraph "##Feedback" immediately followed by a
#      BodyText paragraph with the same contact sentence + hyperlink.
$p18 = $d.Paragraphs.Item(18)
$p18.Range.Text = "##Feedback"
$p18.Style = "First Paragraph"

$p19 = $d.Paragraphs.Item(19)
$p19.Style = "Body Text"

# --- Bookmarks -----------------------------------------------------------
# The source document wraps the three sections in bookmarks
# (clear-screen-and-desk, clear-desk, feedback); the published version
# drops them. Remove any bookmarks the object model exposes.
for ($i = $d.Bookmarks.Count; $i -ge 1; $i--) {
    $d.Bookmarks.Item($i).Delete()
}
